$d = $word.ActiveDocument

# "Versi" + "on" -> merge into a single run reading "Version"
# (also collapses the two runs Word had split the word across).
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Version", 2)

# " 2" (the version number) -> " 1." -- stop short of the trailing "."
# run so we don't swallow the _GoBack bookmark that sits between them.
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " 1.", 2)

# The old trailing "." run (now redundant, since " 1." already ends the
# sentence) sits right after the _GoBack bookmark; delete just that bit
# of text, leaving the bookmark itself and the paragraph mark intact.
$bm = $d.Bookmarks("_GoBack")
$p = $d.Paragraphs(1)
$tail = $d.Range($bm.End, $p.Range.End - 1)
if ($tail.Start -lt $tail.End) {
    $tail.Delete()
}
